# Applies the "dados complementares adicionados no codigo" edit:
#   - Row 20 (AtividadesComplementares): mark column E with a note that the
#     view data was added in code, and flip column F ("AJUSTE EM") from
#     "NÃO" to "SIM".
#   - Move the active selection to A20 (where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("E20").Value = "Dados da View adicionados no código"
$ws.Range("F20").Value = "SIM"

$ws.Range("A20").Select()
